$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 25 (Outdoor plants' locations / creators list update),
# pushing the existing rows 25+ down by one.
$ws.Rows("25:25").Insert() | Out-Null

# Populate the newly inserted row with the new creator entry.
$ws.Range("A25").Value = "Biologist"
$ws.Range("B25").Value = "Saleh Ahmed Saleh"

# Update the active selection / scroll position as recorded in the sheet view.
$ws.Range("C27").Select() | Out-Null
